# Fixed Category duplication bug
#
# Table2 (the Input / Exists / Does Not Exist lookup table in H2:J7) was
# missing a row: "Exercise: Bike, Exercise: Run" -> Exists: "Exercise, Bike, Run".
# Insert that row as the new row 6 of the table, pushing the existing
# "[Blank]" and "Swim" rows down by one, then grow the table to H2:J8.
#
# NOTE: Range.Insert(xlShiftDown) shifts the whole sheet row (all columns),
# which would wrongly disturb Table1 in columns A:F. Table1 must stay put,
# so the row is "inserted" by manually shifting only columns H:J down one
# row at a time (bottom-up, to avoid clobbering) instead of using Insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing Table2 rows (H7:J7 -> H8:J8, H6:J6 -> H7:J7) down by one,
# working bottom-up so we don't overwrite data before it's copied.
$ws.Range("H8").Value2 = $ws.Range("H7").Value2
$ws.Range("I8").Value2 = $ws.Range("I7").Value2
$ws.Range("J8").Value2 = $ws.Range("J7").Value2

$ws.Range("H7").Value2 = $ws.Range("H6").Value2
$ws.Range("I7").Value2 = $ws.Range("I6").Value2
$ws.Range("J7").Value2 = $ws.Range("J6").Value2

# Populate the newly-opened row with the missing lookup entry.
$ws.Range("H6").Value2 = "Exercise: Bike, Exercise: Run"
$ws.Range("J6").Value2 = "Exercise, Bike, Run"

# Grow Table2's range/autofilter to include the new row.
$t2 = $ws.ListObjects.Item("Table2")
$t2.Resize($ws.Range("H2:J8"))

# Match the author's final selection.
$ws.Range("K8").Select()
